# 自动更新Excel文件 - 2025-11-17 23:12:11
# 规则：每行"剩余"天数(E列)递减1；当"剩余"减到1时，说明周期已用完，
# 则以"总天"(D列)重置"剩余"，并把"开始时间"(F列)更新为新的开始日期 20251118。
# 若"开始时间"不是合法的 8 位日期（例如脏数据行），则跳过该行，保持不变。

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStartDate = 20251118

for ($row = 2; $row -le 99; $row++) {
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $total = $dCell.Value2
    $remaining = $eCell.Value2
    $startDate = $fCell.Value2

    if ($null -eq $total -or $null -eq $remaining -or $null -eq $startDate) {
        continue
    }

    # 跳过开始时间不是标准 8 位 YYYYMMDD 格式的脏数据行（例如 202510929）
    $startDateStr = [string]$startDate
    if ($startDateStr.Length -ne 8) {
        continue
    }

    if ($remaining -eq 1) {
        $eCell.Value = $total
        $fCell.Value = $newStartDate
    } else {
        $eCell.Value = $remaining - 1
    }
}
